$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered) from an existing header cell (AC1) onto the
# three new header cells before setting their text, so they match the rest of row 1.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Add new header cells for the season record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record values (Wins=62, Losses=100, Ties=0) for every data row (2-57)
for ($r = 2; $r -le 57; $r++) {
    $ws.Cells.Item($r, 30).Value = 62   # AD
    $ws.Cells.Item($r, 31).Value = 100  # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
